$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "30 Oct 2025, 11:21 AM"

# --- 1 Month Performance sheet: update % Change values and reorder swapped stock names ---
$perf = $wb.Worksheets.Item("1 Month Performance")
$perf.Cells.Item(4, 3).Value = 78.30200000000001
$perf.Cells.Item(5, 3).Value = 66.5852
$perf.Cells.Item(6, 3).Value = 59.6222
$perf.Cells.Item(7, 3).Value = 57.9199
$perf.Cells.Item(9, 3).Value = 54.0542
$perf.Cells.Item(10, 3).Value = 43.1763
$perf.Cells.Item(13, 3).Value = 40.095
$perf.Cells.Item(14, 3).Value = 38.8105
$perf.Cells.Item(15, 2).Value = "TVSSRICHAK"
$perf.Cells.Item(15, 3).Value = 37.5252
$perf.Cells.Item(16, 2).Value = "BHARATSE"
$perf.Cells.Item(16, 3).Value = 37.4206
$perf.Cells.Item(17, 3).Value = 36.7451
$perf.Cells.Item(19, 3).Value = 35.2284
$perf.Cells.Item(20, 3).Value = 33.8004
$perf.Cells.Item(21, 2).Value = "RAMAPHO"
$perf.Cells.Item(21, 3).Value = 33.6546
$perf.Cells.Item(22, 2).Value = "MEGASOFT"
$perf.Cells.Item(22, 3).Value = 33.5822
$perf.Cells.Item(23, 3).Value = 32.6524
$perf.Cells.Item(25, 3).Value = 30.922
$perf.Cells.Item(28, 3).Value = 28.7073
$perf.Cells.Item(30, 2).Value = "EMKAY"
$perf.Cells.Item(30, 3).Value = 26.7268
$perf.Cells.Item(31, 2).Value = "TARACHAND"
$perf.Cells.Item(31, 3).Value = 26.6649
$perf.Cells.Item(32, 2).Value = "ADANIPOWER"
$perf.Cells.Item(32, 3).Value = 26.6572
$perf.Cells.Item(33, 2).Value = "SAGILITY"
$perf.Cells.Item(33, 3).Value = 26.5594
$perf.Cells.Item(34, 3).Value = 26.3781
$perf.Cells.Item(35, 2).Value = "HATSUN"
$perf.Cells.Item(35, 3).Value = 26.0972
$perf.Cells.Item(36, 2).Value = "IFBIND"
$perf.Cells.Item(36, 3).Value = 25.7598
$perf.Cells.Item(37, 2).Value = "INDRAMEDCO"
$perf.Cells.Item(37, 3).Value = 25
$perf.Cells.Item(38, 3).Value = 24.7666
$perf.Cells.Item(39, 2).Value = "SKYGOLD"
$perf.Cells.Item(39, 3).Value = 24.711
$perf.Cells.Item(40, 2).Value = "AUBANK"
$perf.Cells.Item(40, 3).Value = 23.9384
$perf.Cells.Item(41, 2).Value = "CARTRADE"
$perf.Cells.Item(41, 3).Value = 23.8362
$perf.Cells.Item(43, 3).Value = 23.7706
$perf.Cells.Item(44, 3).Value = 23.7317
$perf.Cells.Item(45, 2).Value = "UNIPARTS"
$perf.Cells.Item(45, 3).Value = 23.6096
$perf.Cells.Item(46, 2).Value = "CPEDU"
$perf.Cells.Item(46, 3).Value = 23.5238
$perf.Cells.Item(47, 3).Value = 23.3475
$perf.Cells.Item(48, 3).Value = 23.1174
$perf.Cells.Item(49, 2).Value = "DCBBANK"
$perf.Cells.Item(49, 3).Value = 22.8105
$perf.Cells.Item(50, 2).Value = "ATHERENERG"
$perf.Cells.Item(50, 3).Value = 22.3441
$perf.Cells.Item(51, 3).Value = 22.3209
$perf.Cells.Item(53, 3).Value = 22.029
$perf.Cells.Item(55, 3).Value = 21.6985
$perf.Cells.Item(56, 3).Value = 21.4462
$perf.Cells.Item(58, 3).Value = 21.1934
$perf.Cells.Item(59, 3).Value = 21.1787
$perf.Cells.Item(60, 2).Value = "SHRIRAMFIN"
$perf.Cells.Item(60, 3).Value = 20.8461
$perf.Cells.Item(61, 2).Value = "SKMEGGPROD"
$perf.Cells.Item(61, 3).Value = 20.7435
$perf.Cells.Item(62, 2).Value = "MOLDTECH"
$perf.Cells.Item(62, 3).Value = 20.4604
$perf.Cells.Item(63, 2).Value = "FEDERALBNK"
$perf.Cells.Item(63, 3).Value = 20.4114
$perf.Cells.Item(64, 2).Value = "BANKINDIA"
$perf.Cells.Item(64, 3).Value = 20.2161
$perf.Cells.Item(65, 2).Value = "GRMOVER"
$perf.Cells.Item(65, 3).Value = 20
$perf.Cells.Item(66, 3).Value = 19.8758
$perf.Cells.Item(67, 3).Value = 19.6782
$perf.Cells.Item(68, 3).Value = 19.3989
$perf.Cells.Item(69, 2).Value = "FIVESTAR"
$perf.Cells.Item(69, 3).Value = 19.359
$perf.Cells.Item(70, 2).Value = "REPRO"
$perf.Cells.Item(70, 3).Value = 19.1264
$perf.Cells.Item(71, 2).Value = "BHAGERIA"
$perf.Cells.Item(71, 3).Value = 18.8881
$perf.Cells.Item(72, 2).Value = "ACUTAAS"
$perf.Cells.Item(72, 3).Value = 18.8819
$perf.Cells.Item(73, 2).Value = "WHEELS"
$perf.Cells.Item(73, 3).Value = 18.8605
$perf.Cells.Item(74, 2).Value = "PRECWIRE"
$perf.Cells.Item(74, 3).Value = 18.667

# --- distance from Dma50 sheet: update Distance From Sma50 values ---
$dma = $wb.Worksheets.Item("distance from Dma50")
$dma.Cells.Item(2, 3).Value = 10.1636
$dma.Cells.Item(3, 3).Value = 7.5924
$dma.Cells.Item(4, 3).Value = 6.5153
$dma.Cells.Item(5, 3).Value = 5.3911
$dma.Cells.Item(6, 3).Value = 5.3114
$dma.Cells.Item(7, 3).Value = 5.2126
$dma.Cells.Item(8, 3).Value = 4.5961
$dma.Cells.Item(9, 3).Value = 4.4969
$dma.Cells.Item(10, 3).Value = 3.9517
$dma.Cells.Item(11, 3).Value = 3.6811
$dma.Cells.Item(12, 3).Value = 3.5995
$dma.Cells.Item(13, 3).Value = 3.4744
$dma.Cells.Item(14, 3).Value = 3.2642
$dma.Cells.Item(15, 3).Value = 3.2074
$dma.Cells.Item(16, 3).Value = 3.142
$dma.Cells.Item(17, 3).Value = 2.9738
$dma.Cells.Item(18, 3).Value = 2.8407
$dma.Cells.Item(19, 3).Value = 2.8098
$dma.Cells.Item(20, 3).Value = 2.5079
$dma.Cells.Item(21, 3).Value = 2.3932
$dma.Cells.Item(22, 3).Value = 1.6172
$dma.Cells.Item(23, 3).Value = 1.51
$dma.Cells.Item(24, 3).Value = 1.4478
$dma.Cells.Item(25, 3).Value = 1.198
$dma.Cells.Item(26, 3).Value = 1.1761
$dma.Cells.Item(27, 3).Value = 1.0596
$dma.Cells.Item(28, 3).Value = 0.8689
$dma.Cells.Item(29, 3).Value = 0.4587
$dma.Cells.Item(30, 3).Value = -2.1513
